$wb = $excel.ActiveWorkbook

# Map of row -> new "want to go" count (column F) for the affected rows.
$updates = @{
    2  = 682
    3  = 524
    5  = 22
    6  = 52
    7  = 46
    8  = 3304
    9  = 4257
    10 = 114
}

# The change applies identically to sheet "展览" and sheet "全部类型".
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
